$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5: fix case of tcp/ip to TCP/IP
$ws.Range("D5").Value = "MATLAB의 TCP/IP 통신"

# D37: replace paper title with feed/source name
$ws.Range("D37").Value = "dsba_seminar"

# D45 / E45: update anomaly detection post title and link
$ws.Range("D45").Value = "Anomaly detection - Local Outlier Factor (LOF)"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/106"

# D51 / E51: update matplotlib -> pandas post title and link
$ws.Range("D51").Value = "[pandas] 데이터프레임 컬럼명 또는 행 인덱스 바꾸는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/1323"
